# current_stocks.xlsx - Inventory sheet fixes
# - Layout fix: the original row 2 (a stray placeholder row with A2 = "'")
#   is removed, shifting every subsequent item row up by one.
# - Negative/incorrect stock value fixes applied to several Case/Piece
#   counts after the shift.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the stray first data row (was row 2: A2 = "'"); this shifts all
# rows below it up by one and updates the sheet dimension automatically.
$ws.Rows("2:2").Delete()

# Apply the stock-count corrections (post-shift row numbers).
$ws.Range("B2").Value = 1
$ws.Range("C2").Value = 5

$ws.Range("B3").Value = 1
$ws.Range("C3").Value = 6

$ws.Range("B4").Value = 43

$ws.Range("B10").Value = 320

$ws.Range("B16").Value = 200

$ws.Range("B17").Value = 30

$ws.Range("B20").Value = 25

$ws.Range("B21").Value = 1
